# quarterly.xlsx update: roll the quarterly columns forward by one quarter
# (drop the oldest quarter "فصل چهارم منتهی به 1399/06", add the newest
# quarter "فصل دوم منتهی به 1401/12") and refresh the reported figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$cols = @("E","F","G","H","I","J","K","L","M","N")

# ---------------------------------------------------------------------
# 1. Quarter header labels (row 8 and row 24) shift one column to the
#    left and a new quarter label is appended at column N.
# ---------------------------------------------------------------------
$quarterLabels = @(
    "فصل اول منتهی به 1399/09",
    "فصل دوم منتهی به 1399/12",
    "فصل سوم منتهی به 1400/03",
    "فصل چهارم منتهی به 1400/06",
    "فصل اول منتهی به 1400/09",
    "فصل دوم منتهی به 1400/12",
    "فصل سوم منتهی به 1401/03",
    "فصل چهارم منتهی به 1401/06",
    "فصل اول منتهی به 1401/09",
    "فصل دوم منتهی به 1401/12"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])8").Value = $quarterLabels[$i]
    $ws.Range("$($cols[$i])24").Value = $quarterLabels[$i]
}

# ---------------------------------------------------------------------
# 2. Data rows: values shift one column to the left, newest quarter's
#    figure is appended at column N. Rows not listed (11,13,14,15,18)
#    stay all zero and are unchanged.
# ---------------------------------------------------------------------
$rowData = @{
    10 = @(0, 0, 125961, 129519, 24001, 60072, 23183, 74830, 121837, 57169)
    12 = @(22286, 7037, 6052, 51056, 82, 17319, 5138, 26581, 9021, 6935)
    16 = @(127, 133, 119, 138, 174, 171, 375, 249, 263, 306)
    17 = @(17006, 28103, 15656, 22213, 26261, 17638, 39307, 41923, 29572, 29581)
    19 = @(25328, 52083, 14449, 96811, 53671, 44495, 73283, 36573, 57789, 105708)
    20 = @(64747, 87356, 162237, 299737, 104189, 139695, 141286, 180156, 218482, 199699)
    26 = @(134, 133, 126, 125, 125, 124, 124, 125, 125, 126)
    27 = @(380, 380, 361, 357, 355, 352, 352, 355, 355, 356)
}

foreach ($r in $rowData.Keys) {
    $values = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $values[$i]
    }
}
